# Update cryptos list (GitHub Actions refresh) on worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'73.042.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.60%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'4.003.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'592.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.53%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'160.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.87%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.93%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.26%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'54.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.97%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  +0.77%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'11.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.76%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'4.636.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.06%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "'3.998.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +8.98%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "'14.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.07%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'20.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.33%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "'72.841.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'436.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.03%  "

# Row 22 - PancakeSwap
$ws.Range("E22").Value = "  +13.17%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'96.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.59%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -3.85%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'4.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +19.99%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'14.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.62%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "'11.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "

# Row 28 - LEO
$ws.Range("E28").Value = "  +2.39%  "

# Row 29 - Filecoin
$ws.Range("E29").Value = "  -1.54%  "

# Row 30 - EthereumClassic
$ws.Range("E30").Value = "  +0.21%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.71%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "'13.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.70%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.52%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "'48.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.68%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "'674.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.11%  "

# Row 36 - OKB
$ws.Range("D36").Value = "'70.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.59%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "'0.0₃0881"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.43%  "

# Row 38 - TheGraph
$ws.Range("D38").Value = "'0.440"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -2.55%  "

# Row 40 - ThetaToken
$ws.Range("D40").Value = "'3.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "

# Row 41 - Dai
$ws.Range("E41").Value = "  -0.06%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +3.95%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.23%  "

# Row 44 - THORChain/VeChain swap
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0490"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "

# Row 45 - THORChain/VeChain swap
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D45").Value = "'10.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.14%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  +1.22%  "

# Row 47 - Fetch.AI
$ws.Range("D47").Value = "'2.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.53%  "

# Row 48 - ApeXProtocol
$ws.Range("D48").Value = "'3.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.02%  "

# Row 49 - Stacks
$ws.Range("E49").Value = "  +1.65%  "

# Row 50 - Maker
$ws.Range("D50").Value = "'2.827.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.41%  "

# Row 51 - LidoDAOToken
$ws.Range("E51").Value = "  +4.90%  "
